$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colB = @(0.614130955618208, 0.5715177553169895, 0.5455507417611614, 0.5350193446718663, 0.5332736761901629, 0.5454085068184042, 0.5993972856610696, 0.7068137377750361, 0.7866505503735937, 0.8231651620822902, 0.8370199723928806, 0.8340348787901917, 0.8243044583176413, 0.8183478569877991, 0.7842681252849104, 0.763411109459156, 0.7514332396047223, 0.7473809441194135, 0.7656294614798469, 0.8271617776284472, 0.8675367679372812, 0.8459734710111206, 0.7646265036521527, 0.6775918264195582)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Range("B" + (2 + $i)).Value = $colB[$i]
}

$colC = @(0.1832342419238913, 0.1834208062695879, 0.1835546655792726, 0.1836140920404183, 0.1836242550233358, 0.1835554472460927, 0.1832945753541324, 0.1829352414307941, 0.1827627830519916, 0.1827039420767989, 0.1826844586316412, 0.1826885306277646, 0.1827022832042289, 0.1827110708412079, 0.1827670212153976, 0.1828063523582557, 0.1828308217780616, 0.1828394246247953, 0.1828019744711469, 0.1826981679641122, 0.1826466253885926, 0.1826726502252924, 0.1828039479263097, 0.183016284057743)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Range("C" + (2 + $i)).Value = $colC[$i]
}

$colD = @(0.07841000983695778, 0.07114887349011667, 0.06672343385146462, 0.06492834140335901, 0.06463077046939247, 0.06669919091892496, 0.07589955194617914, 0.09420230747998914, 0.1078092469323906, 0.1140343744685879, 0.116396725575413, 0.1158877279195565, 0.1142286259334213, 0.1132130319579545, 0.1074031263563313, 0.1038479406557968, 0.1018064134179753, 0.1011157590644416, 0.1042260527649859, 0.1147158082063982, 0.1216007714986489, 0.1179234694618714, 0.1040551008143069, 0.08922290566890467)
for ($i = 0; $i -lt $colD.Length; $i++) {
    $ws.Range("D" + (2 + $i)).Value = $colD[$i]
}

$colE = @(0.1181434392460474, 0.1175958507131689, 0.1173192357701041, 0.1172215238498531, 0.1172062061394428, 0.1173178571845028, 0.1179422723128454, 0.119638922931891, 0.1211725591942816, 0.1219324533912491, 0.1222291369464479, 0.1221648440939944, 0.1219566829153003, 0.1218303402539931, 0.1211241488420178, 0.1207068481700126, 0.1204726870272204, 0.1203944110506967, 0.1207506642910054, 0.1220175827629504, 0.1228976171561484, 0.1224231725503842, 0.1207308371088978, 0.1191294609187743)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Range("E" + (2 + $i)).Value = $colE[$i]
}

$colG = @(0.721069912173391, 0.7219813822036016, 0.7229644642199133, 0.7234715481520979, 0.7235621794036646, 0.7229708718471528, 0.7212962973296868, 0.7213734295576018, 0.7234820020392903, 0.7248875193402995, 0.7254839661955117, 0.7253526545871836, 0.7249353024006666, 0.724688024728863, 0.7233991316375779, 0.7227227702057775, 0.7223757528390991, 0.7222654727419524, 0.7227904221747536, 0.7250561460484874, 0.7269111932292844, 0.7258868652027672, 0.7227597064181168, 0.7209924399986361)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Range("G" + (2 + $i)).Value = $colG[$i]
}

$colH = @(0.8202381719016643, 0.8244515666570464, 0.8273655583896513, 0.8286353301942242, 0.8288511477653699, 0.8273823496522894, 0.8216231406174188, 0.8129202179244714, 0.808101568394946, 0.8062506968717855, 0.8055988066149808, 0.8057370248004787, 0.8061960839039841, 0.8064836495342718, 0.8082293867255572, 0.8093876731593639, 0.8100860038832707, 0.8103279643842285, 0.8092610482670324, 0.8060599179198107, 0.8042533432179511, 0.8051914391907786, 0.8093181943997791, 0.8149976416768965)
for ($i = 0; $i -lt $colH.Length; $i++) {
    $ws.Range("H" + (2 + $i)).Value = $colH[$i]
}

$colK = @(0.3183511061105548, 0.281972806424136, 0.2596402299757159, 0.2505409868583115, 0.2490301674380646, 0.2595175078195098, 0.305807356428403, 0.3965947334169471, 0.4632873076485851, 0.493622340732486, 0.5051084671908654, 0.5026347817237422, 0.4945673363250762, 0.4896256344313485, 0.4613047191515989, 0.4439294006675993, 0.4339352568015613, 0.4305513732369377, 0.4457790701258659, 0.4969369746443135, 0.5303649995560136, 0.5125246222366968, 0.4449428496046437, 0.3720346843419691)
for ($i = 0; $i -lt $colK.Length; $i++) {
    $ws.Range("K" + (2 + $i)).Value = $colK[$i]
}

$colL = @(0.1914607565041564, 0.1842882237729668, 0.1799849335866099, 0.1782567048641397, 0.1779712701194427, 0.1799615231497427, 0.1889668187909592, 0.2074224845806043, 0.2214655793321612, 0.2279589056332156, 0.2304328046644883, 0.2298993397716487, 0.2281621343386178, 0.2270999988899405, 0.2210433315035942, 0.2173546054111313, 0.2152428410896476, 0.2145295358736519, 0.2177462533188503, 0.2286719868331062, 0.2359000757283383, 0.2320343328630514, 0.2175691614677362, 0.202344664879476)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Range("L" + (2 + $i)).Value = $colL[$i]
}

$colN = @(1.667097230385806, 1.684216953420999, 1.695265583010079, 1.699903038467695, 1.700681241193982, 1.695327578542393, 1.672888649759188, 1.633146400067221, 1.606542914448634, 1.595003361783499, 1.590714525045097, 1.591634602001176, 1.594648894100922, 1.596505777645272, 1.607308374951766, 1.614079527283746, 1.618027057402438, 1.619372714586562, 1.613353246837488, 1.593761326276815, 1.581428688802697, 1.587967660022097, 1.613681427827799, 1.643441853685123)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Range("N" + (2 + $i)).Value = $colN[$i]
}

$colO = @(3.086034324175728, 3.096513499743168, 3.104520111476873, 3.108178404238984, 3.108809754713008, 3.104567846792804, 3.089321287546909, 3.071894626354862, 3.066692587130433, 3.065976498040044, 3.065942581506675, 3.065939335313828, 3.065968953230652, 3.066017990132309, 3.066772588708432, 3.067658156160661, 3.068322853412013, 3.068574590689934, 3.067547809128683, 3.065953815345807, 3.066294896325132, 3.065986346551199, 3.067597212416075, 3.075273995180453)
for ($i = 0; $i -lt $colO.Length; $i++) {
    $ws.Range("O" + (2 + $i)).Value = $colO[$i]
}

